# Template update: rename the "male" indicator columns from the old
# "...Maso" spelling to the corrected "...Masc" spelling, and rename the
# "CCSARI..." columns to "CCIRAG...". The header row layout (columns A:AD)
# is unchanged positionally - only the text of these specific header
# cells changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("G1").Value  = "ETINumMasc"    # was ETINumMaso
$ws.Range("K1").Value  = "ETIDenoMasc"   # was ETIDenoMaso
$ws.Range("N1").Value  = "HospMasc"      # was HospMaso
$ws.Range("Q1").Value  = "UCIMasc"       # was UCIMaso
$ws.Range("T1").Value  = "DefMasc"       # was DefMaso
$ws.Range("W1").Value  = "NeuMasc"       # was NeuMaso
$ws.Range("Y1").Value  = "CCIRAGFem"     # was CCSARIFem
$ws.Range("Z1").Value  = "CCIRAGMasc"    # was CCSARIMaso
$ws.Range("AA1").Value = "CCIRAGST"      # was CCSARIST
$ws.Range("AC1").Value = "VentMasc"      # was VentMaso
